$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.191680431365967
$ws.Range("B1").Value = 3.078534364700317
$ws.Range("C1").Value = 5.371914386749268
$ws.Range("D1").Value = 2.276468992233276
$ws.Range("E1").Value = 1.387879729270935
